$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new row of data (row 11): date, task, hours
$ws.Range("B11").Value = "15/16-07-2015"
$ws.Range("C11").Value = "EndScreen + polish + clean + proto done"
$ws.Range("D11").Value = 8

# Move the active selection to D12, matching the post-edit UI state
$ws.Range("D12").Select()
